$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 13.19694982948556
$ws.Cells.Item(2, 4).Value = 5.018955193107012
$ws.Cells.Item(2, 5).Value = 16.61403878198083
$ws.Cells.Item(2, 6).Value = 24.83750444975048
$ws.Cells.Item(2, 7).Value = 3.63440757364232
$ws.Cells.Item(2, 9).Value = 26.36578563139043
$ws.Cells.Item(2, 11).Value = 9.604211389286807
$ws.Cells.Item(2, 12).Value = 9.134117306044917
$ws.Cells.Item(2, 13).Value = 14.20287480544433
$ws.Cells.Item(2, 14).Value = 20.01139299716783
$ws.Cells.Item(2, 15).Value = 22.22881989655433

$ws.Cells.Item(3, 2).Value = 13.09555754905248
$ws.Cells.Item(3, 4).Value = 4.966380610283188
$ws.Cells.Item(3, 5).Value = 16.64931619423525
$ws.Cells.Item(3, 6).Value = 24.83341257061285
$ws.Cells.Item(3, 7).Value = 3.636150263713091
$ws.Cells.Item(3, 9).Value = 26.46605124668061
$ws.Cells.Item(3, 11).Value = 9.269166079757039
$ws.Cells.Item(3, 12).Value = 9.1212585355496
$ws.Cells.Item(3, 13).Value = 14.18199927643133
$ws.Cells.Item(3, 14).Value = 20.06926173789802
$ws.Cells.Item(3, 15).Value = 22.27406301617595

$ws.Cells.Item(4, 2).Value = 13.03515972027106
$ws.Cells.Item(4, 4).Value = 4.933281430996626
$ws.Cells.Item(4, 5).Value = 16.6722560046838
$ws.Cells.Item(4, 6).Value = 24.83726515853857
$ws.Cells.Item(4, 7).Value = 3.637278127621229
$ws.Cells.Item(4, 9).Value = 26.53148674962478
$ws.Cells.Item(4, 11).Value = 9.055203005754276
$ws.Cells.Item(4, 12).Value = 9.114737649134753
$ws.Cells.Item(4, 13).Value = 14.17110444140048
$ws.Cells.Item(4, 14).Value = 20.10649778878538
$ws.Cells.Item(4, 15).Value = 22.30653928343834

$ws.Cells.Item(5, 2).Value = 13.01103689166918
$ws.Cells.Item(5, 4).Value = 4.919593369946415
$ws.Cells.Item(5, 5).Value = 16.6819266902098
$ws.Cells.Item(5, 6).Value = 24.84043686187678
$ws.Cells.Item(5, 7).Value = 3.637752331450986
$ws.Cells.Item(5, 9).Value = 26.55912693079324
$ws.Cells.Item(5, 11).Value = 8.966015420158014
$ws.Cells.Item(5, 12).Value = 9.11242824386391
$ws.Cells.Item(5, 13).Value = 14.16715196396841
$ws.Cells.Item(5, 14).Value = 20.1221015874313
$ws.Cells.Item(5, 15).Value = 22.32095291295598

$ws.Cells.Item(6, 2).Value = 13.00706154339875
$ws.Cells.Item(6, 4).Value = 4.917308579172717
$ws.Cells.Item(6, 5).Value = 16.68355200435281
$ws.Cells.Item(6, 6).Value = 24.84106027651839
$ws.Cells.Item(6, 7).Value = 3.63783195515338
$ws.Cells.Item(6, 9).Value = 26.56377546111394
$ws.Cells.Item(6, 11).Value = 8.95108776817526
$ws.Cells.Item(6, 12).Value = 9.112065844631831
$ws.Cells.Item(6, 13).Value = 14.16652519129073
$ws.Cells.Item(6, 14).Value = 20.12471858136792
$ws.Cells.Item(6, 15).Value = 22.32341745480833

$ws.Cells.Item(7, 2).Value = 13.03483237887692
$ws.Cells.Item(7, 4).Value = 4.933097630497818
$ws.Cells.Item(7, 5).Value = 16.67238511991767
$ws.Cells.Item(7, 6).Value = 24.83730144704975
$ws.Cells.Item(7, 7).Value = 3.637284463764962
$ws.Cells.Item(7, 9).Value = 26.53185556692961
$ws.Cells.Item(7, 11).Value = 9.054008163783672
$ws.Cells.Item(7, 12).Value = 9.114705092142799
$ws.Cells.Item(7, 13).Value = 14.17104915927751
$ws.Cells.Item(7, 14).Value = 20.10670648518769
$ws.Cells.Item(7, 15).Value = 22.30672889832219

$ws.Cells.Item(8, 2).Value = 13.16161811983395
$ws.Cells.Item(8, 4).Value = 5.001000764131499
$ws.Cells.Item(8, 5).Value = 16.62593740533822
$ws.Cells.Item(8, 6).Value = 24.83477354681614
$ws.Cells.Item(8, 7).Value = 3.63499647423827
$ws.Cells.Item(8, 9).Value = 26.3995542009432
$ws.Cells.Item(8, 11).Value = 9.490442359723986
$ws.Cells.Item(8, 12).Value = 9.129399769120905
$ws.Cells.Item(8, 13).Value = 14.19528016621901
$ws.Cells.Item(8, 14).Value = 20.03099312465716
$ws.Cells.Item(8, 15).Value = 22.2434439115711

$ws.Cells.Item(9, 2).Value = 13.42395005920364
$ws.Cells.Item(9, 4).Value = 5.127417190921816
$ws.Cells.Item(9, 5).Value = 16.54496703463045
$ws.Cells.Item(9, 6).Value = 24.88022220087881
$ws.Cells.Item(9, 7).Value = 3.630966680224001
$ws.Cells.Item(9, 9).Value = 26.17078717036733
$ws.Cells.Item(9, 11).Value = 10.27792187222862
$ws.Cells.Item(9, 12).Value = 9.169013036551664
$ws.Cells.Item(9, 13).Value = 14.25787914843129
$ws.Cells.Item(9, 14).Value = 19.89598677115033
$ws.Cells.Item(9, 15).Value = 22.15667937800327

$ws.Cells.Item(10, 2).Value = 13.62361609878917
$ws.Cells.Item(10, 4).Value = 5.215863124657283
$ws.Cells.Item(10, 5).Value = 16.49159136450055
$ws.Cells.Item(10, 6).Value = 24.94414503127862
$ws.Cells.Item(10, 7).Value = 3.628281719608804
$ws.Cells.Item(10, 9).Value = 26.02134688326899
$ws.Cells.Item(10, 11).Value = 10.8114053319757
$ws.Cells.Item(10, 12).Value = 9.204542060127089
$ws.Cells.Item(10, 13).Value = 14.31281624286094
$ws.Cells.Item(10, 14).Value = 19.80492747612582
$ws.Cells.Item(10, 15).Value = 22.11577767083952

$ws.Cells.Item(11, 2).Value = 13.71564213576808
$ws.Cells.Item(11, 4).Value = 5.255071694363112
$ws.Cells.Item(11, 5).Value = 16.46862578537777
$ws.Cells.Item(11, 6).Value = 24.97978963507353
$ws.Cells.Item(11, 7).Value = 3.627119530471496
$ws.Cells.Item(11, 9).Value = 25.95739555954376
$ws.Cells.Item(11, 11).Value = 11.04373674440557
$ws.Cells.Item(11, 12).Value = 9.222062663417628
$ws.Cells.Item(11, 13).Value = 14.33969253923082
$ws.Cells.Item(11, 14).Value = 19.76525049086176
$ws.Cells.Item(11, 15).Value = 22.10214284452177

$ws.Cells.Item(12, 2).Value = 13.75063551991741
$ws.Cells.Item(12, 4).Value = 5.26976610331241
$ws.Cells.Item(12, 5).Value = 16.46011761956542
$ws.Cells.Item(12, 6).Value = 24.9942240650741
$ws.Cells.Item(12, 7).Value = 3.626687908793365
$ws.Cells.Item(12, 9).Value = 25.93375744779722
$ws.Cells.Item(12, 11).Value = 11.13018596053836
$ws.Cells.Item(12, 12).Value = 9.228889047086817
$ws.Cells.Item(12, 13).Value = 14.35013556674501
$ws.Cells.Item(12, 14).Value = 19.75047571274219
$ws.Cells.Item(12, 15).Value = 22.09769503873589

$ws.Cells.Item(13, 2).Value = 13.74309304468787
$ws.Cells.Item(13, 4).Value = 5.266608290889866
$ws.Cells.Item(13, 5).Value = 16.46194163776756
$ws.Cells.Item(13, 6).Value = 24.99107382955119
$ws.Cells.Item(13, 7).Value = 3.62678048999873
$ws.Cells.Item(13, 9).Value = 25.93882260229137
$ws.Cells.Item(13, 11).Value = 11.11163619923438
$ws.Cells.Item(13, 12).Value = 9.22741039629388
$ws.Cells.Item(13, 13).Value = 14.34787475381506
$ws.Cells.Item(13, 14).Value = 19.75364662237658
$ws.Cells.Item(13, 15).Value = 22.09862113079035

$ws.Cells.Item(14, 2).Value = 13.7185183089499
$ws.Cells.Item(14, 4).Value = 5.256283709938481
$ws.Cells.Item(14, 5).Value = 16.46792204181024
$ws.Cells.Item(14, 6).Value = 24.98095844094777
$ws.Cells.Item(14, 7).Value = 3.627083851093512
$ws.Cells.Item(14, 9).Value = 25.95543923991768
$ws.Cells.Item(14, 11).Value = 11.05087984450035
$ws.Cells.Item(14, 12).Value = 9.222620453274761
$ws.Cells.Item(14, 13).Value = 14.34054640375822
$ws.Cells.Item(14, 14).Value = 19.76402995608451
$ws.Cells.Item(14, 15).Value = 22.10176258130391

$ws.Cells.Item(15, 2).Value = 13.70348367041534
$ws.Cells.Item(15, 4).Value = 5.249939517465096
$ws.Cells.Item(15, 5).Value = 16.47160972559675
$ws.Cells.Item(15, 6).Value = 24.97488420156158
$ws.Cells.Item(15, 7).Value = 3.627270770905429
$ws.Cells.Item(15, 9).Value = 25.96569277810073
$ws.Cells.Item(15, 11).Value = 11.01346449841459
$ws.Cells.Item(15, 12).Value = 9.219711331003989
$ws.Cells.Item(15, 13).Value = 14.33609199291388
$ws.Cells.Item(15, 14).Value = 19.77042257767357
$ws.Cells.Item(15, 15).Value = 22.10377998700087

$ws.Cells.Item(16, 2).Value = 13.61762376049848
$ws.Cells.Item(16, 4).Value = 5.213279634009146
$ws.Cells.Item(16, 5).Value = 16.49311862290035
$ws.Cells.Item(16, 6).Value = 24.94194707272761
$ws.Cells.Item(16, 7).Value = 3.628358859436713
$ws.Cells.Item(16, 9).Value = 26.02560728940565
$ws.Cells.Item(16, 11).Value = 10.79600997496144
$ws.Cells.Item(16, 12).Value = 9.203424075589117
$ws.Cells.Item(16, 13).Value = 14.31109728116436
$ws.Cells.Item(16, 14).Value = 19.80755550815077
$ws.Cells.Item(16, 15).Value = 22.11676882960297

$ws.Cells.Item(17, 2).Value = 13.56523887598572
$ws.Cells.Item(17, 4).Value = 5.19052324932217
$ws.Cells.Item(17, 5).Value = 16.50664998416531
$ws.Cells.Item(17, 6).Value = 24.92341784933583
$ws.Cells.Item(17, 7).Value = 3.629041503342257
$ws.Cells.Item(17, 9).Value = 26.06339458428956
$ws.Cells.Item(17, 11).Value = 10.65992612541241
$ws.Cells.Item(17, 12).Value = 9.193777754500752
$ws.Cells.Item(17, 13).Value = 14.29624281564969
$ws.Cells.Item(17, 14).Value = 19.8307818618856
$ws.Cells.Item(17, 15).Value = 22.12601085911712

$ws.Cells.Item(18, 2).Value = 13.53522276395211
$ws.Cells.Item(18, 4).Value = 5.177338131125575
$ws.Cells.Item(18, 5).Value = 16.51455670925689
$ws.Cells.Item(18, 6).Value = 24.91337885821029
$ws.Cells.Item(18, 7).Value = 3.629439717867527
$ws.Cells.Item(18, 9).Value = 26.08550817651836
$ws.Cells.Item(18, 11).Value = 10.58068189840657
$ws.Cells.Item(18, 12).Value = 9.188357521200526
$ws.Cells.Item(18, 13).Value = 14.28787670491828
$ws.Cells.Item(18, 14).Value = 19.8443054847262
$ws.Cells.Item(18, 15).Value = 22.13179456724443

$ws.Cells.Item(19, 2).Value = 13.5250802744411
$ws.Cells.Item(19, 4).Value = 5.172857516497493
$ws.Cells.Item(19, 5).Value = 16.5172550833553
$ws.Cells.Item(19, 6).Value = 24.91008627283133
$ws.Cells.Item(19, 7).Value = 3.629575505357629
$ws.Cells.Item(19, 9).Value = 26.093060630157
$ws.Cells.Item(19, 11).Value = 10.55368549430796
$ws.Cells.Item(19, 12).Value = 9.186544427092763
$ws.Cells.Item(19, 13).Value = 14.28507478239966
$ws.Cells.Item(19, 14).Value = 19.84891262684163
$ws.Cells.Item(19, 15).Value = 22.13383317961072

$ws.Cells.Item(20, 2).Value = 13.57080370332516
$ws.Cells.Item(20, 4).Value = 5.192955706977612
$ws.Cells.Item(20, 5).Value = 16.50519673563316
$ws.Cells.Item(20, 6).Value = 24.92532635260855
$ws.Cells.Item(20, 7).Value = 3.628968257953793
$ws.Cells.Item(20, 9).Value = 26.05933280657554
$ws.Cells.Item(20, 11).Value = 10.67451343699677
$ws.Cells.Item(20, 12).Value = 9.194791392023614
$ws.Cells.Item(20, 13).Value = 14.29780573798681
$ws.Cells.Item(20, 14).Value = 19.82829236602011
$ws.Cells.Item(20, 15).Value = 22.1249785948952

$ws.Cells.Item(21, 2).Value = 13.72573278144441
$ws.Cells.Item(21, 4).Value = 5.25932048562132
$ws.Cells.Item(21, 5).Value = 16.46616034486077
$ws.Cells.Item(21, 6).Value = 24.98390422280954
$ws.Cells.Item(21, 7).Value = 3.626994516923063
$ws.Cells.Item(21, 9).Value = 25.95054282511569
$ws.Cells.Item(21, 11).Value = 11.06876727140208
$ws.Cells.Item(21, 12).Value = 9.224022202360244
$ws.Cells.Item(21, 13).Value = 14.34269175593133
$ws.Cells.Item(21, 14).Value = 19.76097334200556
$ws.Cells.Item(21, 15).Value = 22.10082044278452

$ws.Cells.Item(22, 2).Value = 13.82782122956349
$ws.Cells.Item(22, 4).Value = 5.301799334260926
$ws.Cells.Item(22, 5).Value = 16.4417456138814
$ws.Cells.Item(22, 6).Value = 25.02764311998782
$ws.Cells.Item(22, 7).Value = 3.625753938145325
$ws.Cells.Item(22, 9).Value = 25.88281609972666
$ws.Cells.Item(22, 11).Value = 11.31750380944726
$ws.Cells.Item(22, 12).Value = 9.244241859357713
$ws.Cells.Item(22, 13).Value = 14.37357249668777
$ws.Cells.Item(22, 14).Value = 19.71843336775485
$ws.Cells.Item(22, 15).Value = 22.08920163872549

$ws.Cells.Item(23, 2).Value = 13.77326738141273
$ws.Cells.Item(23, 4).Value = 5.279211177534098
$ws.Cells.Item(23, 5).Value = 16.45467600679641
$ws.Cells.Item(23, 6).Value = 25.00380251573904
$ws.Cells.Item(23, 7).Value = 3.626411553894048
$ws.Cells.Item(23, 9).Value = 25.91865461771329
$ws.Cells.Item(23, 11).Value = 11.18557774531777
$ws.Cells.Item(23, 12).Value = 9.233349408672877
$ws.Cells.Item(23, 13).Value = 14.35695137384748
$ws.Cells.Item(23, 14).Value = 19.74100480244796
$ws.Cells.Item(23, 15).Value = 22.09502118015137

$ws.Cells.Item(24, 2).Value = 13.56828752930121
$ws.Cells.Item(24, 4).Value = 5.191856310687873
$ws.Cells.Item(24, 5).Value = 16.50585335253039
$ws.Cells.Item(24, 6).Value = 24.92446160606883
$ws.Cells.Item(24, 7).Value = 3.629001354269572
$ws.Cells.Item(24, 9).Value = 26.06116792398695
$ws.Cells.Item(24, 11).Value = 10.66792164940908
$ws.Cells.Item(24, 12).Value = 9.194332735106709
$ws.Cells.Item(24, 13).Value = 14.29709859864388
$ws.Cells.Item(24, 14).Value = 19.82941733596966
$ws.Cells.Item(24, 15).Value = 22.12544381649992

$ws.Cells.Item(25, 2).Value = 13.35166014049171
$ws.Cells.Item(25, 4).Value = 5.09397376330138
$ws.Cells.Item(25, 5).Value = 16.56579437143268
$ws.Cells.Item(25, 6).Value = 24.86254671733787
$ws.Cells.Item(25, 7).Value = 3.63200822164195
$ws.Cells.Item(25, 9).Value = 26.22939761247352
$ws.Cells.Item(25, 11).Value = 10.07258650570863
$ws.Cells.Item(25, 12).Value = 9.157156078381828
$ws.Cells.Item(25, 13).Value = 14.23935441847798
$ws.Cells.Item(25, 14).Value = 19.93107612681066
$ws.Cells.Item(25, 15).Value = 22.17614465161126
